$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Ranking')
$ws.Range('C2').Value = 2
$ws.Range('E2').Value = 6
$ws.Range('F2').Value = 2
$ws.Range('G2').Value = 0.25
$ws.Range('B3').Value = 'DeepAR'
$ws.Range('C3').Value = 1
$ws.Range('E3').Value = 7
$ws.Range('F3').Value = 1
$ws.Range('G3').Value = 0.125
$ws.Range('H3').Value = 0.003212793603783612
$ws.Range('I3').Value = 0.002468646867616643
$ws.Range('B4').Value = 'AV-MCPS'
$ws.Range('C4').Value = 0
$ws.Range('E4').Value = 8
$ws.Range('F4').Value = 0
$ws.Range('G4').Value = 0
$ws.Range('H4').Value = 0.002343166998705213
$ws.Range('I4').Value = 0.001773557864384127
$ws.Range('B5').Value = 'LSPMW'
$ws.Range('C5').Value = 0
$ws.Range('E5').Value = 8
$ws.Range('F5').Value = 0
$ws.Range('G5').Value = 0
$ws.Range('H5').Value = 0.01088643503280608
$ws.Range('I5').Value = 0.01118631828981691
$ws.Range('C6').Value = 0
$ws.Range('D6').Value = 0
$ws.Range('E6').Value = 8
$ws.Range('F6').Value = 0
$ws.Range('G6').Value = 0
$ws.Range('B7').Value = 'MCPS'
$ws.Range('C7').Value = 0
$ws.Range('D7').Value = 0
$ws.Range('E7').Value = 8
$ws.Range('F7').Value = 0
$ws.Range('G7').Value = 0
$ws.Range('H7').Value = 0.002729303042345392
$ws.Range('I7').Value = 0.002276929651970951
$ws.Range('D8').Value = 1
$ws.Range('E8').Value = 7
$ws.Range('F8').Value = -1
$ws.Range('D9').Value = 1
$ws.Range('E9').Value = 7
$ws.Range('F9').Value = -1
$ws.Range('B10').Value = 'EnCQR-LSTM'
$ws.Range('D10').Value = 1
$ws.Range('E10').Value = 7
$ws.Range('F10').Value = -1
$ws.Range('H10').Value = 0.005726608455807441
$ws.Range('I10').Value = 0.005666169253838035

$ws = $wb.Worksheets.Item('Matriz_Superioridad')
$ws.Range('C2').Value = 0
$ws.Range('D2').Value = 0
$ws.Range('G2').Value = 0
$ws.Range('H2').Value = 0
$ws.Range('J2').Value = 0
$ws.Range('B3').Value = 0
$ws.Range('D3').Value = 0
$ws.Range('E3').Value = 0
$ws.Range('B4').Value = 0
$ws.Range('C4').Value = 0
$ws.Range('E4').Value = 0
$ws.Range('F4').Value = 0
$ws.Range('J4').Value = 0
$ws.Range('C5').Value = 0
$ws.Range('D5').Value = 0
$ws.Range('G5').Value = 0
$ws.Range('H5').Value = 0
$ws.Range('I5').Value = 0
$ws.Range('J5').Value = 0
$ws.Range('D6').Value = 0
$ws.Range('G6').Value = 0
$ws.Range('H6').Value = 0
$ws.Range('I6').Value = 0
$ws.Range('J6').Value = 0
$ws.Range('B7').Value = 0
$ws.Range('E7').Value = 0
$ws.Range('F7').Value = 0
$ws.Range('J7').Value = 0
$ws.Range('B8').Value = 0
$ws.Range('E8').Value = 0
$ws.Range('F8').Value = 0
$ws.Range('J8').Value = 0
$ws.Range('E9').Value = 0
$ws.Range('F9').Value = 0
$ws.Range('J9').Value = 0
$ws.Range('B10').Value = 0
$ws.Range('D10').Value = 0
$ws.Range('E10').Value = 0
$ws.Range('F10').Value = 0
$ws.Range('G10').Value = 0
$ws.Range('H10').Value = 0
$ws.Range('I10').Value = 0

$ws = $wb.Worksheets.Item('Matriz_Pvalores')
$ws.Range('C2').Value = 0.002097018750532698
$ws.Range('D2').Value = 0.004743526100567941
$ws.Range('E2').Value = 0.9130504277383529
$ws.Range('F2').Value = 0.8917294934277238
$ws.Range('G2').Value = 0.003172905104933843
$ws.Range('H2').Value = 0.003988342141515222
$ws.Range('I2').Value = 0.0006952802747612896
$ws.Range('J2').Value = 0.01638296781455018
$ws.Range('B3').Value = 0.002097018750532698
$ws.Range('D3').Value = 0.02957323853369753
$ws.Range('E3').Value = 0.01929193598628043
$ws.Range('F3').Value = 0.000823113484367255
$ws.Range('G3').Value = 0.2005357013888658
$ws.Range('H3').Value = 0.3555178985048615
$ws.Range('I3').Value = 0.223000704671779
$ws.Range('J3').Value = 0.0006251772513028442
$ws.Range('B4').Value = 0.004743526100567941
$ws.Range('C4').Value = 0.02957323853369753
$ws.Range('E4').Value = 0.04081448025316758
$ws.Range('F4').Value = 0.002218914254807469
$ws.Range('G4').Value = 0.1195518541354652
$ws.Range('H4').Value = 0.0857381227163625
$ws.Range('I4').Value = 0.3502234558708817
$ws.Range('J4').Value = 0.0327880517627599
$ws.Range('B5').Value = 0.9130504277383529
$ws.Range('C5').Value = 0.01929193598628043
$ws.Range('D5').Value = 0.04081448025316758
$ws.Range('F5').Value = 0.9863562133347268
$ws.Range('G5').Value = 0.02737605755675454
$ws.Range('H5').Value = 0.02861187543977062
$ws.Range('I5').Value = 0.01770966400219565
$ws.Range('J5').Value = 0.0912879164279572
$ws.Range('B6').Value = 0.8917294934277238
$ws.Range('C6').Value = 0.000823113484367255
$ws.Range('D6').Value = 0.002218914254807469
$ws.Range('E6').Value = 0.9863562133347268
$ws.Range('G6').Value = 0.002240944158671043
$ws.Range('H6').Value = 0.002427679187583021
$ws.Range('I6').Value = 0.002559704324455181
$ws.Range('J6').Value = 0.004861864736348087
$ws.Range('B7').Value = 0.003172905104933843
$ws.Range('C7').Value = 0.2005357013888658
$ws.Range('D7').Value = 0.1195518541354652
$ws.Range('E7').Value = 0.02737605755675454
$ws.Range('F7').Value = 0.002240944158671043
$ws.Range('H7').Value = 0.3219462813494158
$ws.Range('I7').Value = 0.5880646422052729
$ws.Range('J7').Value = 0.002932945703924883
$ws.Range('B8').Value = 0.003988342141515222
$ws.Range('C8').Value = 0.3555178985048615
$ws.Range('D8').Value = 0.0857381227163625
$ws.Range('E8').Value = 0.02861187543977062
$ws.Range('F8').Value = 0.002427679187583021
$ws.Range('G8').Value = 0.3219462813494158
$ws.Range('I8').Value = 0.3792712059915497
$ws.Range('J8').Value = 0.003034218204464256
$ws.Range('B9').Value = 0.0006952802747612896
$ws.Range('C9').Value = 0.223000704671779
$ws.Range('D9').Value = 0.3502234558708817
$ws.Range('E9').Value = 0.01770966400219565
$ws.Range('F9').Value = 0.002559704324455181
$ws.Range('G9').Value = 0.5880646422052729
$ws.Range('H9').Value = 0.3792712059915497
$ws.Range('J9').Value = 0.0664544017438875
$ws.Range('B10').Value = 0.01638296781455018
$ws.Range('C10').Value = 0.0006251772513028442
$ws.Range('D10').Value = 0.0327880517627599
$ws.Range('E10').Value = 0.0912879164279572
$ws.Range('F10').Value = 0.004861864736348087
$ws.Range('G10').Value = 0.002932945703924883
$ws.Range('H10').Value = 0.003034218204464256
$ws.Range('I10').Value = 0.0664544017438875

$ws = $wb.Worksheets.Item('Matriz_DM_Original')
$ws.Range('C2').Value = 7.083283308482033
$ws.Range('D2').Value = 5.679458912127445
$ws.Range('E2').Value = -0.1162589752162882
$ws.Range('F2').Value = -0.144992218631827
$ws.Range('G2').Value = 6.338157550730909
$ws.Range('H2').Value = 5.956120129563035
$ws.Range('I2').Value = 9.464221551842021
$ws.Range('J2').Value = 3.981327185483624
$ws.Range('B3').Value = -7.083283308482033
$ws.Range('D3').Value = -3.312931623975139
$ws.Range('E3').Value = -3.788612492765368
$ws.Range('F3').Value = -9.058396733432964
$ws.Range('G3').Value = -1.530940715526004
$ws.Range('H3').Value = -1.043770623842952
$ws.Range('I3').Value = -1.441086070738292
$ws.Range('J3').Value = -9.7283604744742
$ws.Range('B4').Value = -5.679458912127445
$ws.Range('C4').Value = 3.312931623975139
$ws.Range('E4').Value = -2.978126341872117
$ws.Range('F4').Value = -6.977348038010106
$ws.Range('G4').Value = 1.974499871279546
$ws.Range('H4').Value = 2.269889334820201
$ws.Range('I4').Value = 1.056747819358406
$ws.Range('J4').Value = -3.203602407340061
$ws.Range('B5').Value = 0.1162589752162882
$ws.Range('C5').Value = 3.788612492765368
$ws.Range('D5').Value = 2.978126341872117
$ws.Range('F5').Value = 0.01819296994154352
$ws.Range('G5').Value = 3.396079502860692
$ws.Range('H5').Value = 3.348380570779105
$ws.Range('I5').Value = 3.888729763883627
$ws.Range('J5').Value = 2.21332194743845
$ws.Range('B6').Value = 0.144992218631827
$ws.Range('C6').Value = 9.058396733432964
$ws.Range('D6').Value = 6.977348038010106
$ws.Range('E6').Value = -0.01819296994154352
$ws.Range('G6').Value = 6.958968355485228
$ws.Range('H6').Value = 6.811612805996516
$ws.Range('I6').Value = 6.715615282037597
$ws.Range('J6').Value = 5.641016868496634
$ws.Range('B7').Value = -6.338157550730909
$ws.Range('C7').Value = 1.530940715526004
$ws.Range('D7').Value = -1.974499871279546
$ws.Range('E7').Value = -3.396079502860692
$ws.Range('F7').Value = -6.958968355485228
$ws.Range('H7').Value = 1.129199339886732
$ws.Range('I7').Value = -0.5881017630580124
$ws.Range('J7').Value = -6.474208885539672
$ws.Range('B8').Value = -5.956120129563035
$ws.Range('C8').Value = 1.043770623842952
$ws.Range('D8').Value = -2.269889334820201
$ws.Range('E8').Value = -3.348380570779105
$ws.Range('F8').Value = -6.811612805996516
$ws.Range('G8').Value = -1.129199339886732
$ws.Range('I8').Value = -0.9875688298325775
$ws.Range('J8').Value = -6.415177338567293
$ws.Range('B9').Value = -9.464221551842021
$ws.Range('C9').Value = 1.441086070738292
$ws.Range('D9').Value = -1.056747819358406
$ws.Range('E9').Value = -3.888729763883627
$ws.Range('F9').Value = -6.715615282037597
$ws.Range('G9').Value = 0.5881017630580124
$ws.Range('H9').Value = 0.9875688298325775
$ws.Range('J9').Value = -2.504386448566041
$ws.Range('B10').Value = -3.981327185483624
$ws.Range('C10').Value = 9.7283604744742
$ws.Range('D10').Value = 3.203602407340061
$ws.Range('E10').Value = -2.21332194743845
$ws.Range('F10').Value = -5.641016868496634
$ws.Range('G10').Value = 6.474208885539672
$ws.Range('H10').Value = 6.415177338567293
$ws.Range('I10').Value = 2.504386448566041
